$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "243.83"
$ws.Range("D2").ClearFormats()

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "25.18"
$ws.Range("D3").ClearFormats()

# Row 4
$ws.Range("B4").Value = "HuobiToken"
$ws.Range("C4").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.197"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "3HuobiTokenHT"

# Row 5
$ws.Range("B5").Value = "Cronos"
$ws.Range("C5").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05740"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "4CronosCRO"

# Row 6
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.489"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "5KuCoinTokenKCS"

# Row 7
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.112"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "6GateTokenGT"

# Row 8
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8091"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "7MXTokenMX"

# Row 9
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8394"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "8FTXTokenFTT"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1339"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "9WazirXWRX"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06962"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02840"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "11BitrueCoinBTR"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09360"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "12BitMartTokenBMX"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001519"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "13BitForexTokenBF"

# Row 15
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0005980"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "14OneONE"

# Row 16
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006223"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "15TigerCashTCH"

# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.500"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "16LEOLEO"

# Row 18
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.092"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "17BTSETokenBTSE"

# Row 19
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3196"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "18BitpandaEcosystemTokenBEST"

# Row 20
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03135"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "19LiechtensteinCryptoassetsExchangeLCX"

# Row 21
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1301"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "20ProBitTokenPROB"

# Row 22
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.740"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "21MCDexMCB"

# Row 23
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04655"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "22CoinExTokenCET"

# Row 24
$ws.Range("B24").Value = "ZBToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1329"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "23ZBTokenZB"

# Row 25
$ws.Range("B25").Value = "BitKan"
$ws.Range("C25").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001236"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "24BitKanKAN"

# Row 26
$ws.Range("B26").Value = "HotbitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004262"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "25HotbitTokenHTB"

# Row 27
$ws.Range("B27").Value = "NitroEx"
$ws.Range("C27").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009701"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "26NitroExNTX"

# Row 28
$ws.Range("B28").Value = "UpBots"
$ws.Range("C28").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001501"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "27UpBotsUBXT"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03613"
$ws.Range("D40").ClearFormats()

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006291"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1047"
$ws.Range("D42").ClearFormats()

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002925"
$ws.Range("D43").ClearFormats()

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005281"
$ws.Range("D45").ClearFormats()

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.2500"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002286"
$ws.Range("D48").ClearFormats()
